# Update the "Marking"/"Total" row of the concise marksheet:
# - B11 ("Marking" -> Right marks per correct answer): 3 -> 5
# - B12 ("Total" -> Right marks total): 51 -> 85
# - E12 ("Total" -> correct/total marks display): 42/84 -> 85/140

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 85
$ws.Range("E12").Value = "85/140"
